$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Consolidate the card's separate fields (rows 2-8) into a single
# Python-tuple-style string in A2, then remove the now-empty rows.
$ws.Range("A2").Value = "('Earl of Squirrel', ['{4}{G}{G}', 'Creature — Squirrel Noble Advisor', 'Squirrellink (Damage dealt by this creature also causes you to create that many 1/1 green Squirrel creature tokens.)', 'Creature tokens you control are Squirrels in addition to their other creature types.', 'Other Squirrels you control get +1/+1.', '4/4'])"

$ws.Rows("3:8").Delete()
